# Applies the cryptos-list price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.014.36"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "1.596.07"
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'211.73"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("E7").Value = "  +0.75%  "
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("E9").Value = "  +0.15%  "
$ws.Range("D10").Value = "'18.25"
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("D11").Value = "'0.0808"
$ws.Range("E11").Value = "  +2.16%  "
$ws.Range("D12").Value = "1.820.42"
$ws.Range("E12").Value = "  +0.89%  "
$ws.Range("D13").Value = "1.597.37"
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("E14").Value = "  -0.21%  "
$ws.Range("E15").Value = "  +2.04%  "
$ws.Range("D16").Value = "26.010.87"
$ws.Range("E16").Value = "  +0.49%  "
$ws.Range("D17").Value = "'60.82"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("D18").Value = "0.0₃0729"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("E19").Value = "  -0.11%  "
$ws.Range("D20").Value = "'203.73"
$ws.Range("E20").Value = "  +5.63%  "
$ws.Range("E21").Value = "  +1.84%  "
$ws.Range("D22").Value = "'9.26"
$ws.Range("E22").Value = "  -1.05%  "
$ws.Range("D23").Value = "'6.04"
$ws.Range("E23").Value = "  +1.79%  "
$ws.Range("E24").Value = "  +13.22%  "
$ws.Range("D25").Value = "'144.05"
$ws.Range("E25").Value = "  +1.80%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  -7.40%  "
$ws.Range("D28").Value = "'15.20"
$ws.Range("E28").Value = "  +0.72%  "
$ws.Range("D29").Value = "'6.53"
$ws.Range("E29").Value = "  +1.45%  "
$ws.Range("E30").Value = "  +0.96%  "
$ws.Range("D31").Value = "'0.0477"
$ws.Range("E31").Value = "  +1.15%  "
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("E33").Value = "  -3.70%  "
$ws.Range("E34").Value = "  -0.32%  "
$ws.Range("D36").Value = "1.130.29"
$ws.Range("D37").Value = "'0.0164"
$ws.Range("E37").Value = "  +8.09%  "
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("E39").Value = "  +2.30%  "
$ws.Range("E40").Value = "  -0.88%  "
$ws.Range("E41").Value = "  -1.48%  "
$ws.Range("D42").Value = "'0.778"
$ws.Range("E42").Value = "  -2.44%  "
$ws.Range("E43").Value = "  +0.96%  "
$ws.Range("D44").Value = "1.733.46"
$ws.Range("E44").Value = "  +0.87%  "
$ws.Range("D45").Value = "'92.22"
$ws.Range("E45").Value = "  -1.10%  "
$ws.Range("D46").Value = "'54.27"
$ws.Range("E46").Value = "  +2.33%  "
$ws.Range("D47").Value = "'1.50"
$ws.Range("E47").Value = "  -0.69%  "
$ws.Range("E48").Value = "  -0.39%  "
$ws.Range("B49").Value = "USDD"
$ws.Range("C49").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D49").Value = "'1.01"
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.406"
$ws.Range("E50").Value = "  -0.25%  "
$ws.Range("D51").Value = "0.0₇0948"
$ws.Range("E51").Value = "  -15.15%  "
